$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. All of these columns (D: Price, E: Volume(1h), G: Hora)
# are stored as plain text in the sheet, so force text format before assigning
# to avoid Excel auto-converting numeric-looking strings (e.g. "14", "285.39")
# into numbers/percentages.
$updates = @{
    'D2' = '285.39'
    'E2' = '0.13%'
    'G2' = '14'
    'D3' = '29.36'
    'E3' = '3.70%'
    'G3' = '14'
    'D4' = '5.081'
    'E4' = '1.22%'
    'G4' = '14'
    'D5' = '0.06703'
    'E5' = '2.29%'
    'G5' = '14'
    'E6' = '0.76%'
    'G6' = '14'
    'D7' = '1.411'
    'E7' = '-0.87%'
    'G7' = '14'
    'D8' = '0.8968'
    'E8' = '-3.03%'
    'G8' = '14'
    'D9' = '0.1583'
    'E9' = '1.60%'
    'G9' = '14'
    'D10' = '0.07105'
    'E10' = '9.92%'
    'G10' = '14'
    'D11' = '0.07621'
    'E11' = '0.43%'
    'G11' = '14'
    'D12' = '0.02921'
    'E12' = '1.46%'
    'G12' = '14'
    'D13' = '0.08990'
    'E13' = '0.27%'
    'G13' = '14'
    'D14' = '0.001596'
    'E14' = '0.24%'
    'G14' = '14'
    'D15' = '0.04466'
    'E15' = '1.05%'
    'G15' = '14'
    'D16' = '0.0006484'
    'G16' = '14'
    'D17' = '0.006180'
    'E17' = '2.90%'
    'G17' = '14'
    'D18' = '3.454'
    'E18' = '0.21%'
    'G18' = '14'
    'D19' = '3.442'
    'E19' = '1.96%'
    'G19' = '14'
    'D20' = '2.231'
    'E20' = '-0.37%'
    'G20' = '14'
    'D21' = '0.3232'
    'E21' = '1.15%'
    'G21' = '14'
    'E22' = '1.23%'
    'G22' = '14'
    'D23' = '3.902'
    'E23' = '-1.30%'
    'G23' = '14'
    'D24' = '0.1559'
    'E24' = '1.39%'
    'G24' = '14'
    'D25' = '0.001205'
    'E25' = '1.94%'
    'G25' = '14'
    'D26' = '0.004365'
    'E26' = '-1.83%'
    'G26' = '14'
    'D27' = '0.0001171'
    'E27' = '-5.99%'
    'G27' = '14'
    'E28' = '0.32%'
    'G28' = '14'
    'G29' = '14'
    'G30' = '14'
    'G31' = '14'
    'G32' = '14'
    'G33' = '14'
    'G34' = '14'
    'G35' = '14'
    'G36' = '14'
    'G37' = '14'
    'G38' = '14'
    'G39' = '14'
    'D40' = '0.04247'
    'E40' = '2.34%'
    'G40' = '14'
    'D41' = '0.006779'
    'E41' = '1.36%'
    'G41' = '14'
    'D42' = '0.1236'
    'E42' = '0.17%'
    'G42' = '14'
    'D43' = '0.002231'
    'E43' = '2.28%'
    'G43' = '14'
    'D44' = '0.01273'
    'E44' = '10.85%'
    'G44' = '14'
    'D45' = '0.00005539'
    'E45' = '-1.13%'
    'G45' = '14'
    'D46' = '1.974'
    'E46' = '0.40%'
    'G46' = '14'
    'E47' = '15.88%'
    'G47' = '14'
    'G48' = '14'
    'G49' = '14'
    'G50' = '14'
    'G51' = '14'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}
